# Updates the "cryptos" price/volume table (and swaps the Aave / BabyDogeCoin
# rows) to reflect the latest scrape, per the GitHub Actions commit.
#
# Column D ("Price") holds numeric-looking text (e.g. "1.002", "0.5430") that
# Excel would otherwise silently coerce into real numbers (dropping trailing
# zeros, losing exotic "0.0<sub>5</sub>8385"-style digits, etc.). Prefixing
# those assignments with a leading apostrophe forces Excel to keep them as
# plain text, exactly like a user typing '0.5430 into the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.929.98'
$ws.Range("E2").Value = '  +6.44%  '
$ws.Range("D3").Value = '1.732.14'
$ws.Range("E3").Value = '  +4.54%  '
$ws.Range("D4").Value = "'" + '1.002'
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = "'" + '228.95'
$ws.Range("E5").Value = '  +4.17%  '
$ws.Range("D6").Value = "'" + '0.5430'
$ws.Range("E6").Value = '  +3.69%  '
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("D8").Value = "'" + '0.2772'
$ws.Range("E8").Value = '  +3.62%  '
$ws.Range("D9").Value = "'" + '0.06725'
$ws.Range("E9").Value = '  +5.58%  '
$ws.Range("E10").Value = '  +4.96%  '
$ws.Range("D11").Value = "'" + '0.07835'
$ws.Range("E11").Value = '  +1.48%  '
$ws.Range("D12").Value = "'" + '4.711'
$ws.Range("E12").Value = '  +2.16%  '
$ws.Range("D13").Value = '1.769.13'
$ws.Range("E13").Value = '  +6.78%  '
$ws.Range("D14").Value = '1.969.31'
$ws.Range("E14").Value = '  +4.48%  '
$ws.Range("D15").Value = "'" + '0.6003'
$ws.Range("E15").Value = '  +6.10%  '
$ws.Range("D16").Value = "'" + '0.0' + [char]0x2085 + '8385'
$ws.Range("E16").Value = '  +1.61%  '
$ws.Range("D17").Value = "'" + '68.86'
$ws.Range("E17").Value = '  +5.18%  '
$ws.Range("D18").Value = '27.890.94'
$ws.Range("E18").Value = '  +6.32%  '
$ws.Range("D19").Value = "'" + '216.25'
$ws.Range("E19").Value = '  +12.46%  '
$ws.Range("D20").Value = "'" + '4.825'
$ws.Range("E20").Value = '  +2.49%  '
$ws.Range("E21").Value = '  -0.23%  '
$ws.Range("D22").Value = "'" + '10.92'
$ws.Range("E22").Value = '  +4.72%  '
$ws.Range("D23").Value = "'" + '6.241'
$ws.Range("E23").Value = '  +3.72%  '
$ws.Range("E24").Value = '  -0.18%  '
$ws.Range("D25").Value = "'" + '146.31'
$ws.Range("E25").Value = '  +2.13%  '
$ws.Range("D26").Value = "'" + '0.1246'
$ws.Range("E26").Value = '  +3.71%  '
$ws.Range("D27").Value = "'" + '7.441'
$ws.Range("E27").Value = '  +2.06%  '
$ws.Range("D28").Value = "'" + '1.638'
$ws.Range("E28").Value = '  +8.51%  '
$ws.Range("D29").Value = "'" + '16.88'
$ws.Range("E29").Value = '  +5.62%  '
$ws.Range("D30").Value = "'" + '0.05611'
$ws.Range("E30").Value = '  -0.60%  '
$ws.Range("D31").Value = "'" + '1.317'
$ws.Range("E31").Value = '  +2.95%  '
$ws.Range("D32").Value = "'" + '3.730'
$ws.Range("E32").Value = '  +6.31%  '
$ws.Range("D33").Value = "'" + '3.533'
$ws.Range("E33").Value = '  +5.24%  '
$ws.Range("E34").Value = '  +3.75%  '
$ws.Range("D35").Value = "'" + '0.9836'
$ws.Range("E35").Value = '  +3.96%  '
$ws.Range("E36").Value = '  +1.72%  '
$ws.Range("E37").Value = '  +1.27%  '
$ws.Range("D38").Value = "'" + '0.5912'
$ws.Range("E38").Value = '  +2.25%  '
$ws.Range("D39").Value = "'" + '0.01669'
$ws.Range("E39").Value = '  +4.19%  '
$ws.Range("D40").Value = "'" + '5.942'
$ws.Range("E40").Value = '  +0.27%  '
$ws.Range("D41").Value = '1.044.99'
$ws.Range("E41").Value = '  +2.14%  '
$ws.Range("D42").Value = "'" + '0.8429'
$ws.Range("E42").Value = '  -0.43%  '
$ws.Range("D43").Value = "'" + '1.001'
$ws.Range("E43").Value = '  -0.25%  '
$ws.Range("D44").Value = "'" + '102.52'
$ws.Range("E44").Value = '  +1.03%  '
$ws.Range("D45").Value = '1.874.50'
$ws.Range("E45").Value = '  +4.40%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = "'" + '0.0' + [char]0x2088 + '118'
$ws.Range("E46").Value = '  +12.73%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = "'" + '59.85'
$ws.Range("E47").Value = '  +2.33%  '
$ws.Range("D48").Value = "'" + '8.284'
$ws.Range("E48").Value = '  +3.03%  '
$ws.Range("D49").Value = "'" + '0.4416'
$ws.Range("E49").Value = '  +1.52%  '
$ws.Range("D50").Value = "'" + '1.003'
$ws.Range("E50").Value = '  -0.35%  '
$ws.Range("D51").Value = "'" + '0.05315'
$ws.Range("E51").Value = '  -0.16%  '
